$d = $word.ActiveDocument

# The document carries a "_GoBack" bookmark (Word's automatic "last edit
# location" marker) that currently sits near the end of the "Madame
# Vandevorst" paragraph (right after "+ diagramme de classes "). Re-saving
# the document after further edits elsewhere moved that marker back to the
# very start of the document (inside the first paragraph, right before the
# logo picture). Reproduce that relocation here: delete the old bookmark and
# re-insert a collapsed "_GoBack" bookmark at the top of the document.

$old = $d.Bookmarks("_GoBack")
$old.Delete()

$r = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $r)
